# Add a "success" column (D) that flags the first 9 "list" entries
# (list index 1 through 8) as successes ("1") and the rest ("0"),
# including the list-index-0 row which is not a success.
#
# success = 1 if 1 <= list_index <= 8 else 0, for list_index = A value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header - same bold/centered/bordered look as the other header cells
$ws.Range("D1").Value = "success"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows: sheet row 2 => list index 0 ... sheet row 25 => list index 23
for ($row = 2; $row -le 25; $row++) {
    $listIndex = $row - 2
    if ($listIndex -ge 1 -and $listIndex -le 8) {
        $flag = "1"
    } else {
        $flag = "0"
    }

    # Write through a text formula, then paste-special as values, so the
    # numeric-looking "0"/"1" lands as literal text (shared string) instead
    # of being auto-coerced to a number - matching the source data's dtype.
    $cell = $ws.Cells.Item($row, 4)
    $cell.Formula = "=""" + $flag + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
